$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new values ---
$ws.Range("A2").Value = 45114.50694444445
$ws.Range("B2").Value = 24.021
$ws.Range("C2").Value = 16.687
$ws.Range("D2").Value = 4.255
$ws.Range("E2").Value = 50.696
$ws.Range("F2").Value = 41.946
$ws.Range("G2").Value = 18.904
$ws.Range("H2").Value = 63.181
$ws.Range("I2").Value = 29.086
$ws.Range("J2").Value = 12.432
$ws.Range("K2").Value = 19.179
$ws.Range("L2").Value = 19.925
$ws.Range("M2").Value = 20.923
$ws.Range("N2").Value = 6.036
$ws.Range("O2").Value = 18.798
$ws.Range("P2").Value = 26.508
$ws.Range("Q2").Value = 15.598
$ws.Range("R2").Value = 3.832
$ws.Range("S2").Value = 2.606
$ws.Range("T2").Value = 278.79
$ws.Range("U2").Value = 52.368
$ws.Range("V2").Value = 17.351
$ws.Range("W2").Value = 34.892
$ws.Range("X2").Value = 18.126
$ws.Range("Y2").Value = 2.393
$ws.Range("Z2").Value = 31.289
$ws.Range("AA2").Value = 15.326
$ws.Range("AB2").Value = 13.706
$ws.Range("AC2").Value = 16.029
$ws.Range("AD2").Value = 20.711
$ws.Range("AE2").Value = 3.641
$ws.Range("AF2").Value = 55.941
$ws.Range("AG2").Value = 9.720000000000001
$ws.Range("AH2").Value = 21.693

$ws.Range("A3").Value = 45114.51388888889
$ws.Range("B3").Value = 13.452
$ws.Range("C3").Value = 9.433999999999999
$ws.Range("D3").Value = 1.716
$ws.Range("E3").Value = 28.652
$ws.Range("F3").Value = 23.745
$ws.Range("G3").Value = 10.586
$ws.Range("H3").Value = 43.392
$ws.Range("I3").Value = 16.288
$ws.Range("J3").Value = 7.014
$ws.Range("K3").Value = 10.695
$ws.Range("L3").Value = 11.379
$ws.Range("M3").Value = 11.97
$ws.Range("N3").Value = 3.383
$ws.Range("O3").Value = 10.527
$ws.Range("P3").Value = 14.853
$ws.Range("Q3").Value = 8.986000000000001
$ws.Range("R3").Value = 1.613
$ws.Range("S3").Value = 0.978
$ws.Range("T3").Value = 152.921
$ws.Range("U3").Value = 29.534
$ws.Range("V3").Value = 9.717000000000001
$ws.Range("W3").Value = 19.595
$ws.Range("X3").Value = 10.45
$ws.Range("Y3").Value = 1.316
$ws.Range("Z3").Value = 20.448
$ws.Range("AA3").Value = 8.583
$ws.Range("AB3").Value = 7.78
$ws.Range("AC3").Value = 9.099
$ws.Range("AD3").Value = 11.916
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 39.253
$ws.Range("AG3").Value = 5.402
$ws.Range("AH3").Value = 12.149

$ws.Range("A4").Value = 45114.52083333334
$ws.Range("B4").Value = 10.569
$ws.Range("C4").Value = 7.495
$ws.Range("D4").Value = 1.127
$ws.Range("E4").Value = 22.594
$ws.Range("F4").Value = 18.733
$ws.Range("G4").Value = 8.318
$ws.Range("H4").Value = 33.954
$ws.Range("I4").Value = 12.798
$ws.Range("J4").Value = 5.545
$ws.Range("K4").Value = 8.407
$ws.Range("L4").Value = 9.022
$ws.Range("M4").Value = 9.481
$ws.Range("N4").Value = 2.658
$ws.Range("O4").Value = 8.271000000000001
$ws.Range("P4").Value = 11.676
$ws.Range("Q4").Value = 7.1
$ws.Range("R4").Value = 1.041
$ws.Range("S4").Value = 0.639
$ws.Range("T4").Value = 118.573
$ws.Range("U4").Value = 23.179
$ws.Range("V4").Value = 7.635
$ws.Range("W4").Value = 15.39
$ws.Range("X4").Value = 8.275
$ws.Range("Y4").Value = 1.026
$ws.Range("Z4").Value = 15.949
$ws.Range("AA4").Value = 6.744
$ws.Range("AB4").Value = 6.103
$ws.Range("AC4").Value = 7.144
$ws.Range("AD4").Value = 9.477
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 30.612
$ws.Range("AG4").Value = 4.239
$ws.Range("AH4").Value = 9.545

$ws.Range("A5").Value = 45114.52777777778
$ws.Range("B5").Value = 4.8
$ws.Range("C5").Value = 3.28
$ws.Range("D5").Value = 0.72
$ws.Range("E5").Value = 10.15
$ws.Range("F5").Value = 8.44
$ws.Range("G5").Value = 3.78
$ws.Range("H5").Value = 17.93
$ws.Range("I5").Value = 5.82
$ws.Range("J5").Value = 2.49
$ws.Range("K5").Value = 3.76
$ws.Range("L5").Value = 4.06
$ws.Range("M5").Value = 4.23
$ws.Range("N5").Value = 1.21
$ws.Range("O5").Value = 3.76
$ws.Range("P5").Value = 5.29
$ws.Range("Q5").Value = 3.33
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.37
$ws.Range("T5").Value = 49.92
$ws.Range("U5").Value = 10.68
$ws.Range("V5").Value = 3.47
$ws.Range("W5").Value = 7.01
$ws.Range("X5").Value = 3.83
$ws.Range("Y5").Value = 0.44
$ws.Range("Z5").Value = 8.15
$ws.Range("AA5").Value = 3.07
$ws.Range("AB5").Value = 2.84
$ws.Range("AC5").Value = 3.31
$ws.Range("AD5").Value = 4.27
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 16.3
$ws.Range("AG5").Value = 1.88
$ws.Range("AH5").Value = 4.34

# --- Remove row 6 (only 4 data rows remain) ---
$ws.Rows(6).Delete()

# --- Column width adjustments (AutoFit-equivalent widths) ---
$ws.Columns("B").ColumnWidth = 7.166666666666667
$ws.Columns("C").ColumnWidth = 7.166666666666667
$ws.Columns("G").ColumnWidth = 7.166666666666667
$ws.Columns("I").ColumnWidth = 7.166666666666667
$ws.Columns("J").ColumnWidth = 7.166666666666667
$ws.Columns("K").ColumnWidth = 7.166666666666667
$ws.Columns("L").ColumnWidth = 7.166666666666667
$ws.Columns("M").ColumnWidth = 7.166666666666667
$ws.Columns("O").ColumnWidth = 7.166666666666667
$ws.Columns("P").ColumnWidth = 7.166666666666667
$ws.Columns("Q").ColumnWidth = 7.166666666666667
$ws.Columns("T").ColumnWidth = 8.166666666666666
$ws.Columns("V").ColumnWidth = 7.166666666666667
$ws.Columns("X").ColumnWidth = 7.166666666666667
$ws.Columns("Z").ColumnWidth = 7.166666666666667
$ws.Columns("AA").ColumnWidth = 7.166666666666667
$ws.Columns("AB").ColumnWidth = 7.166666666666667
$ws.Columns("AC").ColumnWidth = 7.166666666666667
$ws.Columns("AD").ColumnWidth = 7.166666666666667
$ws.Columns("AH").ColumnWidth = 7.166666666666667
